# The source workbook "A 33655-2023.xlsx" had its data rows (2-13) reshuffled:
# each row 2-13 in the edited file now carries the values that a different
# row used to carry (rows 8 and 9 are unchanged). This mirrors an upstream
# re-sort/re-sync of the export where the row identity (the "Id" in column A
# and everything else in that row) moved to a new position.
#
# Mapping of: destination row -> source row (values to place there)
#   2  <- 12
#   3  <- 2
#   4  <- 3
#   5  <- 11
#   6  <- 5
#   7  <- 10
#   8  <- 8   (unchanged, left alone)
#   9  <- 9   (unchanged, left alone)
#   10 <- 13
#   11 <- 7
#   12 <- 4
#   13 <- 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row. Rows 8 and 9 are not part of the shuffle
# (they keep their original data) so they are intentionally omitted here -
# touching them would risk subtly changing cell typing (e.g. a numeric-
# looking text value silently turning into a real number on read/write)
# for no reason.
$mapping = @{
    2  = 12
    3  = 2
    4  = 3
    5  = 11
    6  = 5
    7  = 10
    10 = 13
    11 = 7
    12 = 4
    13 = 6
}

# Columns Y ("Startdatum") and AA ("Slutdatum") hold a plain date-string
# ("2023-09-11") that, in this sheet, is identical on every single data row.
# Copying row-wide arrays through .Value2 re-parses any text that looks like
# a date into a date serial number (the same way typing a date into a fresh
# cell would), which would needlessly change both the stored type and
# on-disk representation of those two columns. Since the value never
# actually differs between source and destination row here, the simplest
# robust fix is to just not touch columns Y/AA at all - every row keeps the
# date string it already had.
$segments = @(
    @{ first = "A";  last = "X" },   # up to (not including) Y
    @{ first = "Z";  last = "Z" },   # between Y and AA
    @{ first = "AB"; last = "AY" }   # after AA, to the end
)

# Snapshot every source row's segments as arrays BEFORE writing anything, so
# overlapping cycles in the permutation don't clobber values still needed
# (e.g. row 2's original data feeds row 3, row 3's feeds row 4, etc.).
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($seg in $segments) {
            $addr = "$($seg.first)$srcRow`:$($seg.last)$srcRow"
            $rowData[$seg.first] = $ws.Range($addr).Value2
        }
        $snapshot[$srcRow] = $rowData
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($seg in $segments) {
        $addr = "$($seg.first)$destRow`:$($seg.last)$destRow"
        $ws.Range($addr).Value2 = $rowData[$seg.first]
    }
}
